$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1131.4667
$ws.Range("I32").Value = 1239.6
$ws.Range("J32").Value = 1077.4
$ws.Range("K32").Value = 1239.6
$ws.Range("L32").Value = 1077.4
$ws.Range("M32").Value = -913.5999999999999
$ws.Range("N32").Value = -1729.4
$ws.Range("H40").Value = 1871.0416
$ws.Range("I40").Value = 2706.6667
$ws.Range("J40").Value = 1369.6666
$ws.Range("K40").Value = 2706.6667
$ws.Range("L40").Value = 1369.6666
$ws.Range("M40").Value = -2531.6667
$ws.Range("N40").Value = -1719.6666
$ws.Range("H43").Value = 10566.333
$ws.Range("I43").Value = 4600
$ws.Range("J43").Value = 14828
$ws.Range("K43").Value = 4600
$ws.Range("L43").Value = 14828
$ws.Range("M43").Value = -4531
$ws.Range("N43").Value = -14966
$ws.Range("H118").Value = 2291.6453
$ws.Range("I118").Value = 910
$ws.Range("J118").Value = 3164.2632
$ws.Range("K118").Value = 2730
$ws.Range("L118").Value = 9492.7896
$ws.Range("M118").Value = -1073
$ws.Range("N118").Value = -12806.7896
$ws.Range("H138").Value = 8335492.5
$ws.Range("I138").Value = 1898.1428
$ws.Range("J138").Value = 15627387
$ws.Range("K138").Value = 5694.428400000001
$ws.Range("L138").Value = 46882161
$ws.Range("M138").Value = -554.4284000000007
$ws.Range("N138").Value = -46892441

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1618.5
$ws.Range("I110").Value = 977.75
$ws.Range("K110").Value = 977.75
$ws.Range("M110").Value = 1067.25
$ws.Range("H132").Value = 8067161.5
$ws.Range("I132").Value = 9617003
$ws.Range("J132").Value = 7982.4
$ws.Range("K132").Value = 28851009
$ws.Range("L132").Value = 23947.2
$ws.Range("M132").Value = -28848479
$ws.Range("N132").Value = -29007.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11497333
$ws.Range("I31").Value = 2486.7917
$ws.Range("J31").Value = 66672596
$ws.Range("K31").Value = 2486.7917
$ws.Range("L31").Value = 66672596
$ws.Range("M31").Value = -2191.7917
$ws.Range("N31").Value = -66673186
$ws.Range("H34").Value = 11497333
$ws.Range("I34").Value = 2486.7917
$ws.Range("J34").Value = 66672596
$ws.Range("K34").Value = 2486.7917
$ws.Range("L34").Value = 66672596
$ws.Range("M34").Value = -2284.7917
$ws.Range("N34").Value = -66673000
$ws.Range("H58").Value = 1972.1111
$ws.Range("I58").Value = 978.8570999999999
$ws.Range("J58").Value = 5448.5
$ws.Range("K58").Value = 978.8570999999999
$ws.Range("L58").Value = 5448.5
$ws.Range("M58").Value = -775.8570999999999
$ws.Range("N58").Value = -5854.5
$ws.Range("H107").Value = 476.13333
$ws.Range("I107").Value = 358.57144
$ws.Range("J107").Value = 579
$ws.Range("K107").Value = 358.57144
$ws.Range("L107").Value = 579
$ws.Range("M107").Value = 1561.42856
$ws.Range("N107").Value = -4419
$ws.Range("H132").Value = 2691.3872
$ws.Range("I132").Value = 1938.84
$ws.Range("K132").Value = 5816.52
$ws.Range("M132").Value = -3286.52
$ws.Range("H136").Value = 1972.1111
$ws.Range("I136").Value = 978.8570999999999
$ws.Range("J136").Value = 5448.5
$ws.Range("K136").Value = 2936.5713
$ws.Range("L136").Value = 16345.5
$ws.Range("M136").Value = -386.5712999999996
$ws.Range("N136").Value = -21445.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5565.231
$ws.Range("I3").Value = 4626.875
$ws.Range("J3").Value = 7066.6
$ws.Range("K3").Value = 13880.625
$ws.Range("L3").Value = 21199.8
$ws.Range("M3").Value = -13768.625
$ws.Range("N3").Value = -21423.8
$ws.Range("H5").Value = 710.6667
$ws.Range("I5").Value = 566
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 1698
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = -1586
$ws.Range("N5").Value = -3224
$ws.Range("H118").Value = 1830.95
$ws.Range("I118").Value = 329
$ws.Range("J118").Value = 2096
$ws.Range("K118").Value = 987
$ws.Range("L118").Value = 6288
$ws.Range("M118").Value = 256
$ws.Range("N118").Value = -8774
$ws.Range("H122").Value = 719
$ws.Range("I122").Value = 545
$ws.Range("J122").Value = 788.6
$ws.Range("K122").Value = 4905
$ws.Range("L122").Value = 7097.400000000001
$ws.Range("M122").Value = -2455
$ws.Range("N122").Value = -11997.4
$ws.Range("H127").Value = 783
$ws.Range("J127").Value = 783
$ws.Range("L127").Value = 2349
$ws.Range("N127").Value = -12269
$ws.Range("H133").Value = 153850180
$ws.Range("I133").Value = 222225390
$ws.Range("J133").Value = 5975
$ws.Range("K133").Value = 666676170
$ws.Range("L133").Value = 17925
$ws.Range("M133").Value = -666671110
$ws.Range("N133").Value = -28045
$ws.Range("H134").Value = 4004
$ws.Range("I134").Value = 2712.5
$ws.Range("K134").Value = 8137.5
$ws.Range("M134").Value = -3067.5
$ws.Range("H135").Value = 710.6667
$ws.Range("I135").Value = 566
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 5094
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -2559
$ws.Range("N135").Value = -14070
$ws.Range("H136").Value = 3483.8064
$ws.Range("I136").Value = 1269.0769
$ws.Range("J136").Value = 5083.3335
$ws.Range("K136").Value = 3807.2307
$ws.Range("L136").Value = 15250.0005
$ws.Range("M136").Value = 1292.7693
$ws.Range("N136").Value = -25450.0005
$ws.Range("H137").Value = 6948887.5
$ws.Range("J137").Value = 4957.0713
$ws.Range("L137").Value = 14871.2139
$ws.Range("N137").Value = -25071.2139
$ws.Range("H138").Value = 2685
$ws.Range("I138").Value = 1452
$ws.Range("J138").Value = 8850
$ws.Range("K138").Value = 4356
$ws.Range("L138").Value = 26550
$ws.Range("M138").Value = 784
$ws.Range("N138").Value = -36830
$ws.Range("H139").Value = 1707.75
$ws.Range("I139").Value = 1183.238
$ws.Range("J139").Value = 2709.0908
$ws.Range("K139").Value = 3549.714
$ws.Range("L139").Value = 8127.2724
$ws.Range("M139").Value = 1590.286
$ws.Range("N139").Value = -18407.2724

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2008.4615
$ws.Range("I16").Value = 2052.1738
$ws.Range("J16").Value = 1673.3334
$ws.Range("K16").Value = 2052.1738
$ws.Range("L16").Value = 1673.3334
$ws.Range("M16").Value = -1882.1738
$ws.Range("N16").Value = -2013.3334
$ws.Range("H46").Value = 1100
$ws.Range("I46").Value = 845
$ws.Range("J46").Value = 1304
$ws.Range("K46").Value = 845
$ws.Range("L46").Value = 1304
$ws.Range("M46").Value = -657
$ws.Range("N46").Value = -1680
$ws.Range("H82").Value = 2103.75
$ws.Range("I82").Value = 1753
$ws.Range("J82").Value = 2279.125
$ws.Range("K82").Value = 1753
$ws.Range("L82").Value = 2279.125
$ws.Range("M82").Value = -1392
$ws.Range("N82").Value = -3001.125
$ws.Range("H85").Value = 2103.75
$ws.Range("I85").Value = 1753
$ws.Range("J85").Value = 2279.125
$ws.Range("K85").Value = 1753
$ws.Range("L85").Value = 2279.125
$ws.Range("M85").Value = -505
$ws.Range("N85").Value = -4775.125
$ws.Range("H136").Value = 31263656
$ws.Range("I136").Value = 55560868
$ws.Range("J136").Value = 24384.143
$ws.Range("K136").Value = 166682604
$ws.Range("L136").Value = 73152.429
$ws.Range("M136").Value = -166680054
$ws.Range("N136").Value = -78252.429
